$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Fgfr1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = [double]"0.300794"
$ws.Cells.Item(2, 8).Value = [double]"0.902382"
$ws.Cells.Item(2, 9).Value = [double]"0.003836667885433928"
$ws.Cells.Item(2, 10).Value = [double]"0.003836667885433928"
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = [double]"3.76917"
$ws.Cells.Item(2, 14).Value = [double]"11.30751"
$ws.Cells.Item(2, 15).Value = [double]"0.02686645020528053"
$ws.Cells.Item(2, 16).Value = [double]"0.02686645020528053"
$ws.Cells.Item(2, 17).Value = [double]"1.13374372098"
$ws.Cells.Item(2, 18).Value = [double]"10.20369348882"
$ws.Cells.Item(2, 19).Value = [double]"0.0001030776466982096"
$ws.Cells.Item(2, 20).Value = [double]"0.0001030776466982096"

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Fgfr1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = [double]"0.300794"
$ws.Cells.Item(3, 8).Value = [double]"0.902382"
$ws.Cells.Item(3, 9).Value = [double]"0.003836667885433928"
$ws.Cells.Item(3, 10).Value = [double]"0.003836667885433928"
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = [double]"82.48060333333333"
$ws.Cells.Item(3, 14).Value = [double]"247.44181"
$ws.Cells.Item(3, 15).Value = [double]"0.5879175050094569"
$ws.Cells.Item(3, 16).Value = [double]"0.587917505009457"
$ws.Cells.Item(3, 17).Value = [double]"24.80967059904667"
$ws.Cells.Item(3, 18).Value = [double]"223.28703539142"
$ws.Cells.Item(3, 19).Value = [double]"0.002255644210754224"
$ws.Cells.Item(3, 20).Value = [double]"0.002255644210754224"

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Fgfr1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = [double]"0.300794"
$ws.Cells.Item(4, 8).Value = [double]"0.902382"
$ws.Cells.Item(4, 9).Value = [double]"0.003836667885433928"
$ws.Cells.Item(4, 10).Value = [double]"0.003836667885433928"
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = [double]"0.5716056666666667"
$ws.Cells.Item(4, 14).Value = [double]"1.714817"
$ws.Cells.Item(4, 15).Value = [double]"0.004074375838860061"
$ws.Cells.Item(4, 16).Value = [double]"0.004074375838860062"
$ws.Cells.Item(4, 17).Value = [double]"0.1719355548993333"
$ws.Cells.Item(4, 18).Value = [double]"1.547419994094"
$ws.Cells.Item(4, 19).Value = [double]"1.563202693414232e-05"
$ws.Cells.Item(4, 20).Value = [double]"1.563202693414232e-05"

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Fgfr1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = [double]"0.300794"
$ws.Cells.Item(5, 8).Value = [double]"0.902382"
$ws.Cells.Item(5, 9).Value = [double]"0.003836667885433928"
$ws.Cells.Item(5, 10).Value = [double]"0.003836667885433928"
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = [double]"53.47143866666666"
$ws.Cells.Item(5, 14).Value = [double]"160.414316"
$ws.Cells.Item(5, 15).Value = [double]"0.3811416689464024"
$ws.Cells.Item(5, 16).Value = [double]"0.3811416689464024"
$ws.Cells.Item(5, 17).Value = [double]"16.08388792230133"
$ws.Cells.Item(5, 18).Value = [double]"144.754991300712"
$ws.Cells.Item(5, 19).Value = [double]"0.001462314001047352"
$ws.Cells.Item(5, 20).Value = [double]"0.001462314001047352"

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Fgfr1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = [double]"6.268658666666667"
$ws.Cells.Item(6, 8).Value = [double]"18.805976"
$ws.Cells.Item(6, 9).Value = [double]"0.07995758356598558"
$ws.Cells.Item(6, 10).Value = [double]"0.07995758356598558"
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = [double]"3.76917"
$ws.Cells.Item(6, 14).Value = [double]"11.30751"
$ws.Cells.Item(6, 15).Value = [double]"0.02686645020528053"
$ws.Cells.Item(6, 16).Value = [double]"0.02686645020528053"
$ws.Cells.Item(6, 17).Value = [double]"23.62764018664"
$ws.Cells.Item(6, 18).Value = [double]"212.64876167976"
$ws.Cells.Item(6, 19).Value = [double]"0.002148176437410108"
$ws.Cells.Item(6, 20).Value = [double]"0.002148176437410108"

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Fgfr1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = [double]"6.268658666666667"
$ws.Cells.Item(7, 8).Value = [double]"18.805976"
$ws.Cells.Item(7, 9).Value = [double]"0.07995758356598558"
$ws.Cells.Item(7, 10).Value = [double]"0.07995758356598558"
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = [double]"82.48060333333333"
$ws.Cells.Item(7, 14).Value = [double]"247.44181"
$ws.Cells.Item(7, 15).Value = [double]"0.5879175050094569"
$ws.Cells.Item(7, 16).Value = [double]"0.587917505009457"
$ws.Cells.Item(7, 17).Value = [double]"517.0427489173956"
$ws.Cells.Item(7, 18).Value = [double]"4653.384740256561"
$ws.Cells.Item(7, 19).Value = [double]"0.04700846303669939"
$ws.Cells.Item(7, 20).Value = [double]"0.04700846303669941"

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Fgfr1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = [double]"6.268658666666667"
$ws.Cells.Item(8, 8).Value = [double]"18.805976"
$ws.Cells.Item(8, 9).Value = [double]"0.07995758356598558"
$ws.Cells.Item(8, 10).Value = [double]"0.07995758356598558"
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = [double]"0.5716056666666667"
$ws.Cells.Item(8, 14).Value = [double]"1.714817"
$ws.Cells.Item(8, 15).Value = [double]"0.004074375838860061"
$ws.Cells.Item(8, 16).Value = [double]"0.004074375838860062"
$ws.Cells.Item(8, 17).Value = [double]"3.583200816265778"
$ws.Cells.Item(8, 18).Value = [double]"32.24880734639201"
$ws.Cells.Item(8, 19).Value = [double]"0.0003257772466148859"
$ws.Cells.Item(8, 20).Value = [double]"0.000325777246614886"

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Fgfr1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = [double]"6.268658666666667"
$ws.Cells.Item(9, 8).Value = [double]"18.805976"
$ws.Cells.Item(9, 9).Value = [double]"0.07995758356598558"
$ws.Cells.Item(9, 10).Value = [double]"0.07995758356598558"
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = [double]"53.47143866666666"
$ws.Cells.Item(9, 14).Value = [double]"160.414316"
$ws.Cells.Item(9, 15).Value = [double]"0.3811416689464024"
$ws.Cells.Item(9, 16).Value = [double]"0.3811416689464024"
$ws.Cells.Item(9, 17).Value = [double]"335.1941974169351"
$ws.Cells.Item(9, 18).Value = [double]"3016.747776752416"
$ws.Cells.Item(9, 19).Value = [double]"0.03047516684526118"
$ws.Cells.Item(9, 20).Value = [double]"0.03047516684526118"

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Fgfr1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(10, 7).Value = [double]"0.04541633333333334"
$ws.Cells.Item(10, 8).Value = [double]"0.136249"
$ws.Cells.Item(10, 9).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(10, 10).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = [double]"3.76917"
$ws.Cells.Item(10, 14).Value = [double]"11.30751"
$ws.Cells.Item(10, 15).Value = [double]"0.02686645020528053"
$ws.Cells.Item(10, 16).Value = [double]"0.02686645020528053"
$ws.Cells.Item(10, 17).Value = [double]"0.17118188111"
$ws.Cells.Item(10, 18).Value = [double]"1.54063692999"
$ws.Cells.Item(10, 19).Value = [double]"1.556350446372418e-05"
$ws.Cells.Item(10, 20).Value = [double]"1.556350446372418e-05"

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Ncam1"
$ws.Cells.Item(11, 3).Value = "Fgfr1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(11, 7).Value = [double]"0.04541633333333334"
$ws.Cells.Item(11, 8).Value = [double]"0.136249"
$ws.Cells.Item(11, 9).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(11, 10).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = [double]"82.48060333333333"
$ws.Cells.Item(11, 14).Value = [double]"247.44181"
$ws.Cells.Item(11, 15).Value = [double]"0.5879175050094569"
$ws.Cells.Item(11, 16).Value = [double]"0.587917505009457"
$ws.Cells.Item(11, 17).Value = [double]"3.745966574521111"
$ws.Cells.Item(11, 18).Value = [double]"33.71369917069001"
$ws.Cells.Item(11, 19).Value = [double]"0.0003405755745028738"
$ws.Cells.Item(11, 20).Value = [double]"0.0003405755745028738"

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Ncam1"
$ws.Cells.Item(12, 3).Value = "Fgfr1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(12, 7).Value = [double]"0.04541633333333334"
$ws.Cells.Item(12, 8).Value = [double]"0.136249"
$ws.Cells.Item(12, 9).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(12, 10).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = [double]"0.5716056666666667"
$ws.Cells.Item(12, 14).Value = [double]"1.714817"
$ws.Cells.Item(12, 15).Value = [double]"0.004074375838860061"
$ws.Cells.Item(12, 16).Value = [double]"0.004074375838860062"
$ws.Cells.Item(12, 17).Value = [double]"0.02596023349255556"
$ws.Cells.Item(12, 18).Value = [double]"0.233642101433"
$ws.Cells.Item(12, 19).Value = [double]"2.360251022017235e-06"
$ws.Cells.Item(12, 20).Value = [double]"2.360251022017236e-06"

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Ncam1"
$ws.Cells.Item(13, 3).Value = "Fgfr1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 7).Value = [double]"0.04541633333333334"
$ws.Cells.Item(13, 8).Value = [double]"0.136249"
$ws.Cells.Item(13, 9).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(13, 10).Value = [double]"0.0005792914339187697"
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = [double]"53.47143866666666"
$ws.Cells.Item(13, 14).Value = [double]"160.414316"
$ws.Cells.Item(13, 15).Value = [double]"0.3811416689464024"
$ws.Cells.Item(13, 16).Value = [double]"0.3811416689464024"
$ws.Cells.Item(13, 17).Value = [double]"2.428476682298222"
$ws.Cells.Item(13, 18).Value = [double]"21.856290140684"
$ws.Cells.Item(13, 19).Value = [double]"0.0002207921039301545"
$ws.Cells.Item(13, 20).Value = [double]"0.0002207921039301545"

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Ncam1"
$ws.Cells.Item(14, 3).Value = "Fgfr1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = [double]"71.78493233333333"
$ws.Cells.Item(14, 8).Value = [double]"215.354797"
$ws.Cells.Item(14, 9).Value = [double]"0.9156264571146617"
$ws.Cells.Item(14, 10).Value = [double]"0.9156264571146617"
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = [double]"3.76917"
$ws.Cells.Item(14, 14).Value = [double]"11.30751"
$ws.Cells.Item(14, 15).Value = [double]"0.02686645020528053"
$ws.Cells.Item(14, 16).Value = [double]"0.02686645020528053"
$ws.Cells.Item(14, 17).Value = [double]"270.56961340283"
$ws.Cells.Item(14, 18).Value = [double]"2435.12652062547"
$ws.Cells.Item(14, 19).Value = [double]"0.02459963261670848"
$ws.Cells.Item(14, 20).Value = [double]"0.02459963261670849"

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Ncam1"
$ws.Cells.Item(15, 3).Value = "Fgfr1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = [double]"71.78493233333333"
$ws.Cells.Item(15, 8).Value = [double]"215.354797"
$ws.Cells.Item(15, 9).Value = [double]"0.9156264571146617"
$ws.Cells.Item(15, 10).Value = [double]"0.9156264571146617"
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = [double]"82.48060333333333"
$ws.Cells.Item(15, 14).Value = [double]"247.44181"
$ws.Cells.Item(15, 15).Value = [double]"0.5879175050094569"
$ws.Cells.Item(15, 16).Value = [double]"0.587917505009457"
$ws.Cells.Item(15, 17).Value = [double]"5920.864529095841"
$ws.Cells.Item(15, 18).Value = [double]"53287.78076186257"
$ws.Cells.Item(15, 19).Value = [double]"0.5383128221875004"
$ws.Cells.Item(15, 20).Value = [double]"0.5383128221875005"

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Ncam1"
$ws.Cells.Item(16, 3).Value = "Fgfr1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = [double]"71.78493233333333"
$ws.Cells.Item(16, 8).Value = [double]"215.354797"
$ws.Cells.Item(16, 9).Value = [double]"0.9156264571146617"
$ws.Cells.Item(16, 10).Value = [double]"0.9156264571146617"
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = [double]"0.5716056666666667"
$ws.Cells.Item(16, 14).Value = [double]"1.714817"
$ws.Cells.Item(16, 15).Value = [double]"0.004074375838860061"
$ws.Cells.Item(16, 16).Value = [double]"0.004074375838860062"
$ws.Cells.Item(16, 17).Value = [double]"41.03267410301655"
$ws.Cells.Item(16, 18).Value = [double]"369.294066927149"
$ws.Cells.Item(16, 19).Value = [double]"0.003730606314289016"
$ws.Cells.Item(16, 20).Value = [double]"0.003730606314289016"

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Ncam1"
$ws.Cells.Item(17, 3).Value = "Fgfr1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = [double]"71.78493233333333"
$ws.Cells.Item(17, 8).Value = [double]"215.354797"
$ws.Cells.Item(17, 9).Value = [double]"0.9156264571146617"
$ws.Cells.Item(17, 10).Value = [double]"0.9156264571146617"
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = [double]"53.47143866666666"
$ws.Cells.Item(17, 14).Value = [double]"160.414316"
$ws.Cells.Item(17, 15).Value = [double]"0.3811416689464024"
$ws.Cells.Item(17, 16).Value = [double]"0.3811416689464024"
$ws.Cells.Item(17, 17).Value = [double]"3838.44360645265"
$ws.Cells.Item(17, 18).Value = [double]"34545.99245807384"
$ws.Cells.Item(17, 19).Value = [double]"0.3489833959961637"
$ws.Cells.Item(17, 20).Value = [double]"0.3489833959961637"

